$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Fill in the new key/value pair on row 7 (previously an empty placeholder row)
$ws.Range("A7").Value = "style"
$ws.Range("B7").Value = "default"

# Copy the key-column formatting (bold orange font) from an existing key cell
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# Add a new empty placeholder row below, matching the previous template row
$ws.Range("A8").Value = ""
